$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnsNumberParameters")

# Update B20:B22 with new education estimate values.
# These were previously plain numbers; they become text-typed cells
# (quote-prefixed) holding the new estimate values, matching the
# existing convention used elsewhere on the sheet (e.g. B26, B27).
$ws.Cells.Item(20, 2).Value = "'20"
$ws.Cells.Item(21, 2).Value = "'18"
$ws.Cells.Item(22, 2).Value = "'22"

# Move the active selection to B23 to reflect the author's cursor position.
$ws.Range("B23").Select()
